$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2 gains a literal 0; B2 loses its "=C2-2" formula and becomes a
# plain literal value; C2:E2 get refreshed (recalculated) literal values;
# F2 is unchanged.
$ws.Range("A2").Value = 0
$ws.Range("B2").ClearContents()
$ws.Range("B2").Value = 3.793611978731801
$ws.Range("C2").Value = 5.793611978731801
$ws.Range("D2").Value = 7.793611978731802
$ws.Range("E2").Value = 9.793611978731802
$ws.Range("F2").Value = 11.7936119787318

# Row 3
$ws.Range("B3").Value = 1.367211748755173
$ws.Range("C3").Value = 2.270285792844648
$ws.Range("D3").Value = 3.53482330417652
$ws.Range("E3").Value = 5.20058241653066
$ws.Range("F3").Value = 7.285956045259299

# Row 4
$ws.Range("B4").Value = 1.025137670248806
$ws.Range("C4").Value = 1.755588205080587
$ws.Range("D4").Value = 2.812128913965231
$ws.Range("E4").Value = 4.243822837148998
$ws.Range("F4").Value = 6.079625054756848

# Row 5
$ws.Range("B5").Value = 0.7583933316883401
$ws.Range("C5").Value = 1.339873030477419
$ws.Range("D5").Value = 2.209786504348833
$ws.Range("E5").Value = 3.424417302186311
$ws.Range("F5").Value = 5.022569086173307

# Row 6
$ws.Range("B6").Value = 0.55377445281194
$ws.Range("C6").Value = 1.009286703292296
$ws.Range("D6").Value = 1.714809822031961
$ws.Range("E6").Value = 2.731321721486744
$ws.Range("F6").Value = 4.10608223442945

# Row 7
$ws.Range("B7").Value = 0.3993385411285419
$ws.Range("C7").Value = 0.7505088202816107
$ws.Range("D7").Value = 1.313978510606229
$ws.Range("E7").Value = 2.152710186451342
$ws.Range("F7").Value = 3.320464391256474
